# Updates cryptos list values (Price and Volume(1h) columns) per the
# scraped-data refresh. D-column price cells are plain text in this
# sheet (values such as "27.352.38" or "1.000" aren't valid numbers,
# or would lose their exact formatting if coerced to one), so we force
# text interpretation via a temporary Text number format and then
# clear the format again so the cell's style stays at its original
# (default) index -- only the cell's value/text actually changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.352.38"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -4.01%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.861.48"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -4.92%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -1.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.47"
$ws.Range("D5").ClearFormats()

$ws.Range("E6").Value = "  -0.94%  "

$ws.Range("E7").Value = "  -5.81%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3870"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -4.76%  "

$ws.Range("E9").Value = "  -11.62%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07904"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -7.02%  "

$ws.Range("E11").Value = "  -3.60%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.42"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -4.38%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.854.62"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -5.67%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.898"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -4.46%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.160"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -5.44%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9999"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.32%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001034"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -3.59%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "85.77"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -5.48%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06533"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.60%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.16"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -7.25%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.97%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.530"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -5.50%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.357.43"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -4.09%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.86"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -5.19%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.272"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.96%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.086.30"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -4.96%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "152.59"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.19%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.77"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.61%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.064"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -5.33%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.512"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -6.00%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "120.93"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.84%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.492"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.42%  "

$ws.Range("E33").Value = "  -3.53%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9368"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -5.21%  "

$ws.Range("E35").Value = "  -2.39%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.289"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -6.13%  "

$ws.Range("E37").Value = "  -4.16%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06020"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.24%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.219"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.61%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.296"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -8.97%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9999"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.97%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5912"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -5.13%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1888"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.43%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.17"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -9.32%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.279"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -5.78%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5639"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -5.30%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "11.94"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -8.81%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.927"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -6.61%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.366"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.30%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06799"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.58%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "107.91"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.09%  "
